# docs(sprint): Update kanban board for end of sprint 3
#
# End-of-sprint-3 kanban update on the ARCHIVE board: the items below were
# finished during the sprint, so flip their STATUS (column E) from
# "In Progress" to "Complete". Also leave the sheet scrolled/selected at
# the spot where the work left off (row ~31, cell E37 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("ARCHIVE")

# Rows whose STATUS column (E) flips from "In Progress" to "Complete".
$rows = @(19, 27, 33, 34, 35, 37)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "Complete"
}

# Restore the view state: scrolled so row 31 is at the top, E37 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E37").Select()
